$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 9 data
$ws.Range("A9").Value = $ws.Range("A8").Value()
$ws.Range("B9").Value = "Find First and Last Position of Element in Sorted Array"
$ws.Range("C9").Value = "discrete binary search, l&r, equal_range;"

# Copy the formatting of B7 (same alternating highlight style used for B5/B7) onto B9
$ws.Range("B7").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Adjust column B width to match new content
$ws.Columns("B").ColumnWidth = 46.5

# Update the active selection to C9 like in the diff
$ws.Range("C9").Select()
